# Mise a jour de l'application
# Adds 12 new "Entrainement" rows (J-3, 2025-09-24) for rows 553-564 on Feuil1,
# mirroring the formatting of the existing rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 553
$lastNewRow  = 564
$templateRow = 552

# --- Columns shared by every new row -------------------------------------------------
$ws.Range("A$firstNewRow`:A$lastNewRow").Value = "Entrainement"
$ws.Range("C$firstNewRow`:C$lastNewRow").Value = "Global"

# Date column (B) keeps the same date-only display format as the rows above it.
$ws.Cells.Item($templateRow, 2).Copy()
$ws.Range("B$firstNewRow`:B$lastNewRow").PasteSpecial(-4122)
$ws.Range("B$firstNewRow`:B$lastNewRow").Value = 45924

# "MD" column (D) keeps the centered style used by the rows above it.
$ws.Cells.Item($templateRow, 4).Copy()
$ws.Range("D$firstNewRow`:D$lastNewRow").PasteSpecial(-4122)
$ws.Range("D$firstNewRow`:D$lastNewRow").Value = "J-3"

$ws.Application.CutCopyMode = 0

# --- Per-row data: Nom du joueur, Poste, Temps joue, then the 15 numeric metrics -----
$rows = @(
    ,("Malik Boussaid",     "right back",      "01:20:46", 6.03, 0.38, 5.63, 0.35, 0.05, 0,    0,    1, 4.12, 25.25, 4.78, 72, 14, 58, 28)
    ,("Mattheo Haon",       "right back",      "01:21:46", 5.96, 0.35, 5.61, 0.29, 0.06, 0,    0,    0, 4.31, 24.65, 4.75, 44, 10, 35,  6)
    ,("Kamal Bafounta",     "center midfield", "01:21:00", 5.14, 0.21, 4.92, 0.18, 0.03, 0,    0,    0, 3.75, 22.32, 4.48, 27,  4, 16,  2)
    ,("Omar Benyounes",     "center midfield", "01:22:14", 5.88, 0.33, 5.54, 0.24, 0.09, 0.01, 0,    2, 3.53, 26.69, 5.56, 38, 10, 44, 14)
    ,("Fareh Wael",         "center midfield", "01:22:13", 5.21, 0.32, 4.88, 0.26, 0.07, 0,    0,    0, 3.74, 24.23, 4.95, 32,  6, 18,  7)
    ,("Ilan Ihaddadene",    "center midfield", "01:04:03", 4.23, 0.22, 4.01, 0.21, 0.02, 0,    0,    0, 3.62, 23.03, 5.75, 67, 12, 60, 14)
    ,("Yoann Martelat",     "center midfield", "01:20:26", 5.40, 0.21, 5.18, 0.22, 0,    0,    0,    0, 3.93, 20.25, 4.13, 20,  2, 19,  4)
    ,("Jeremie Laurent",    "left forward",    "01:19:39", 5.21, 0.34, 4.86, 0.30, 0.05, 0,    0,    1, 3.82, 25.60, 5.20, 38, 18, 25, 11)
    ,("Emmanuel Valey",     "left forward",    "01:21:47", 5.69, 0.28, 5.40, 0.27, 0.02, 0,    0,    0, 3.78, 22.91, 5.30, 70, 24, 62, 13)
    ,("Karim Belmahi",      "left forward",    "01:21:26", 5.29, 0.27, 5.01, 0.23, 0.05, 0,    0,    0, 3.81, 22.54, 4.93, 29, 11, 26,  8)
    ,("Hedi Nasri",         "right back",      "01:21:26", 5.00, 0.25, 4.75, 0.21, 0.04, 0,    0,    0, 3.56, 24.00, 4.61, 28,  5, 26,  2)
    ,("Naim Ighbane",       "center back",     "01:21:06", 5.64, 0.17, 5.46, 0.14, 0.03, 0,    0,    0, 3.32, 23.37, 4.65, 27,  5, 25,  6)
)

$numCols = $rows[0].Count
$numRows = $rows.Count
$arr = New-Object 'object[,]' $numRows, $numCols

for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r, $c] = $rows[$r][$c]
    }
}

$target = $ws.Range("E$firstNewRow`:V$lastNewRow")
$target.Value = $arr

# --- Replicate the resulting view/selection state -------------------------------------
[void]$ws.Range("D570").Select()
$excel.ActiveWindow.ScrollRow = 536
